$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns stay formatted as text so values like "1.20" or
# "0.940" are not auto-converted to numbers, matching the original inlineStr cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '57.895.71'
$ws.Range('E2').Value = '  +1.45%  '
$ws.Range('D3').Value = '3.070.77'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '514.83'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').Value = '141.02'
$ws.Range('E6').Value = '  +1.13%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '0.435'
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').Value = '7.32'
$ws.Range('E9').Value = '  +2.02%  '
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('D11').Value = '0.378'
$ws.Range('E11').Value = '  +2.70%  '
$ws.Range('D12').Value = '3.597.78'
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('D14').Value = '26.66'
$ws.Range('E14').Value = '  +5.59%  '
$ws.Range('D15').Value = '0.0000164'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '57.906.66'
$ws.Range('E16').Value = '  +1.28%  '
$ws.Range('D17').Value = '3.076.49'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '6.16'
$ws.Range('E18').Value = '  +3.92%  '
$ws.Range('D19').Value = '12.84'
$ws.Range('E19').Value = '  -2.04%  '
$ws.Range('D20').Value = '8.08'
$ws.Range('E20').Value = '  -0.50%  '
$ws.Range('D21').Value = '332.98'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').Value = '0.501'
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').Value = '65.03'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('E25').Value = '  +2.92%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').Value = '  -2.91%  '
$ws.Range('D28').Value = '6.43'
$ws.Range('E28').Value = '  +0.95%  '
$ws.Range('E29').Value = '  +3.60%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').Value = '1.20'
$ws.Range('E31').Value = '  +2.67%  '
$ws.Range('D32').Value = '20.78'
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('D33').Value = '154.92'
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '5.99'
$ws.Range('E35').Value = '  +2.34%  '
$ws.Range('B36').Value = 'EnergySwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D36').Value = '27.19'
$ws.Range('E36').Value = '  +2.55%  '
$ws.Range('E37').Value = '  +3.07%  '
$ws.Range('E38').Value = '  +0.65%  '
$ws.Range('D39').Value = '3.112.04'
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('D40').Value = '3.89'
$ws.Range('E40').Value = '  +1.49%  '
$ws.Range('D41').Value = '36.55'
$ws.Range('E41').Value = '  -0.96%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = '0.654'
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('D44').Value = '2.278.32'
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('E46').Value = '  +1.34%  '
$ws.Range('D47').Value = '20.49'
$ws.Range('E47').Value = '  +3.37%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = '0.940'
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = '5.94'
$ws.Range('E49').Value = '  +1.58%  '
$ws.Range('E50').Value = '  +6.43%  '
$ws.Range('D51').Value = '256.68'
$ws.Range('E51').Value = '  +9.72%  '
